$d = $word.ActiveDocument

# 1) CONSTRUTORA -> PLAENGE
$d.Content.Find.Execute("CONSTRUTORA", $true, $false, $false, $false, $false, $true, 1, $false, "PLAENGE", 2)

# 2) PROJETO -> MOMENTUM
$d.Content.Find.Execute("PROJETO", $true, $false, $false, $false, $false, $true, 1, $false, "MOMENTUM", 2)
